$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Quantum Computing: Unraveling Possibilities" "The Fascinating Story of flight"

# --- Author paragraph: "Dr" + "." + " Alan Turing" -> "Sarah Johnson" (collapse 3 runs into 1) ---
$rng = $d.Content
$rng.Find.Execute("Dr. Alan Turing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Sarah Johnson"

# --- Email paragraph ---
Replace-Text "alanturing@ai-innovations" "sarahjohnson@amail"
Replace-Text "org" "com"

# --- Body paragraph 1 (introduction) ---
Replace-Text "The realm of quantum computing stands as a captivating frontier in the scientific landscape, holding the promise to revolutionize various disciplines" "Since its inception, humanity has looked to the skies with admiration"
Replace-Text "It challenges conventional computing paradigms by harnessing the extraordinary properties of quantum mechanics, opening doors to unprecedented computational power and transformative applications" "While looking up it was noticed that birds glided effortlessly through the air"
Replace-Text "In this essay, we delve into the captivating world of quantum computing, exploring its fundamental principles, its potential impact across diverse fields, and the challenges that lie ahead" "This led to the desire to harness the same power"

# Insert two new sentences after "...harness the same power."
$rng = $d.Content
$rng.Find.Execute("This led to the desire to harness the same power.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" In the realm of science, aviation emerged as researchers delved into the secrets of flight. From the groundbreaking experiments of the Wright brothers to the supersonic marvels of today, this essay provides an exploration into the wondrous world of flight.")

# First break section: "Quantum computing unveils..." -> insert extra break + "Early attempts..."
$rng = $d.Content
$rng.Find.Execute("Quantum computing unveils a universe where subatomic particles, such as electrons or photons, exist in multiple states simultaneously--a phenomenon known as superposition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertBefore([char]11)
$rng = $d.Content
$rng.Find.Execute("Quantum computing unveils a universe where subatomic particles, such as electrons or photons, exist in multiple states simultaneously--a phenomenon known as superposition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Early attempts at mimicking avian locomotion dates as far back as the Renaissance period"

Replace-Text "This enables quantum systems to process exponentially larger amounts of data concurrently, promising solutions to problems that defy classical computation" "During this time, inventors like Leonardo Da Vinci sketched and conceptualized flying machines, setting the stage for future advancements"
Replace-Text "Additionally, the concept of quantum entanglement allows particles to remain interconnected, even when physically separated, enabling communication and computations beyond classical limits" "By the 19th century, scientists began understanding the principles of aerodynamics, which laid the foundation for the first successful heavier-than-air craft"

# Insert two new sentences after "...heavier-than-air craft."
$rng = $d.Content
$rng.Find.Execute("By the 19th century, scientists began understanding the principles of aerodynamics, which laid the foundation for the first successful heavier-than-air craft.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" The Wright brothers' Kitty Hawk experiment in 1903 signified a pivotal moment in aviation history. They managed to achieve controlled, sustained flight, forever changing the dynamics of travel, warfare, and societal progress.")

# Second break section: "The implications..." -> insert extra break + "With the dawn of the 20th century..."
$rng = $d.Content
$rng.Find.Execute("The implications of quantum computing are as vast as they are intriguing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertBefore([char]11)
$rng = $d.Content
$rng.Find.Execute("The implications of quantum computing are as vast as they are intriguing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "With the dawn of the 20th century, aviation witnessed an unprecedented surge in innovation"

Replace-Text "In the realm of materials science, it could accelerate the design of novel materials with enhanced properties, revolutionizing industries ranging from electronics to medicine" "The advent of the jet engine in the 1930s revolutionized the industry, enabling aircraft to reach remarkable speeds and fly at higher altitudes"
Replace-Text "Drug discovery stands to benefit from quantum simulations that elucidate complex molecular interactions, leading to more effective and personalized treatments" "The visionary minds of designers and engineers brought about iconic aircraft such as the Spitfire and the P-51 Mustang"
Replace-Text "Artificial intelligence algorithms, empowered by quantum enhancements, could achieve unprecedented levels of efficiency and accuracy, driving transformative advancements in fields such as natural language processing and image recognition" "With the passing of time, technological advancements gave rise to commercial aviation, transforming travel by offering speed, convenience, and connectivity across the world"

# --- Summary paragraph ---
Replace-Text "Quantum computing presents a paradigm shift in computational possibilities, leveraging the enigmatic principles of quantum mechanics to transcend the limitations of classical computing" "The evolution of flight stands as an enduring testament to human ingenuity and persistence"
Replace-Text "Its potential impact reverberates across diverse fields, holding the promise of transformative breakthroughs in materials science, drug discovery, artificial intelligence, and cryptography" "From the early dreams of flight to the supersonic marvels of today, aviation has forever changed the fabric of humanity"
Replace-Text "While challenges remain in harnessing and controlling quantum systems, the allure of quantum " "Its profound impact on society, warfare, and exploration is evident in the "
Replace-Text "computing continues to inspire researchers and innovators worldwide" "modern world"
Replace-Text "This enthralling field promises to redefine the boundaries of computation, unlocking a new era of scientific discovery and technological advancement" "The story of flight continues to be written, with the skies holding vast potential for future discoveries and innovations"

# Remove the trailing "Total Word Count..." line (and its two line breaks), then add a blank paragraph.
$rng = $d.Content
$rng.Find.Execute(" The story of flight continues to be written, with the skies holding vast potential for future discoveries and innovations.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastParaEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$delRng = $d.Range($rng.End, $lastParaEnd)
$delRng.Delete()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter() | Out-Null

Write-Output $d.Content.Text
